$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-looking string as TEXT (avoid Excel's automatic
# date serial conversion) by staging it in a scratch cell formatted as
# Text, then pasting VALUES ONLY into the destination (preserves the
# destination's existing cell style/format).
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = 0
    $scratch.ClearContents() | Out-Null
}

# 1. Ativacao date: 01/01/2012 -> 01/01/2023
Set-TextValue $ws.Range("B8") "01/01/2023"
Set-TextValue $ws.Range("C8") "01/01/2023"

# 2. Objetivos docente: Paulo Atsushi Suzuki -> Luiz Tadeu Fernandes Eleno
$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# 3. New row 11 (Objectives: english objectives text), copy formatting from row 10
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B11").Value = "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."
$ws.Range("C11").Value = "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."

# 4. Row 13 (Programa resumido), was "Semestral" -> "01/01/2023"
Set-TextValue $ws.Range("B13") "01/01/2023"
Set-TextValue $ws.Range("C13") "01/01/2023"

# 5. New row 14 (Short syllabus), copy formatting from row 13
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B14").Value = "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."
$ws.Range("C14").Value = "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."

# 6. Row 15 (Programa), was "01/01/2012" -> "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("B15").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C15").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# 7. New row 16 (Syllabus), copy formatting from row 15
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B16").Value = "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."
$ws.Range("C16").Value = "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."

$excel.CutCopyMode = 0
